$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the date string shared by A2:A5 (날짜 column) from 2025-11-29 to 2025-12-01.
# A plain Range("A2:A5").Value = "2025-12-01" assignment gets auto-converted
# by Excel into a date serial number (since the text looks like a date), which
# would change the cell type/format and is not what the source workbook has
# (a plain text shared string). To force it to stay literal text, enter it as
# a formula returning a text literal, then copy/paste-special as values so the
# formula collapses down to a static text value (still text, not re-parsed as
# a date) without touching the cell's number format/style.
$ws.Range("A2:A5").Formula = "=""2025-12-01"""
$ws.Range("A2:A5").Copy() | Out-Null
$ws.Range("A2:A5").PasteSpecial(-4163) | Out-Null

# Update the N column (최종점수) values in rows 2-5
$ws.Range("N2").Value = 85.87246918135976
$ws.Range("N3").Value = 85.87246918135976
$ws.Range("N4").Value = 85.87246918135976
$ws.Range("N5").Value = 85.87246918135976
